# The author re-typed a single letter "i" in the middle of the word
# "avez" (turning "Vous avez supprimé..." into "Vous aviez supprimé...",
# i.e. present -> imperfect tense). Word's automatic "_GoBack" bookmark
# (which always marks the location of the last edit) follows that
# keystroke, moving from the end of the following paragraph to right
# after the newly typed "i".

$d = $word.ActiveDocument

# Locate "Vous av" (the part of "avez" that stays before the insertion
# point) and collapse the range to its end, i.e. right before "ez...".
$r = $d.Content
$r.Find.Execute("Vous av", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)

# Type the missing "i" - InsertAfter behaves like typing: it expands the
# range to cover the freshly inserted text.
$r.InsertAfter("i")

# Mark the boundaries of that new "i" with temporary/real bookmarks so it
# ends up in its own run (matching the way Word splits runs around
# bookmarks), then drop the temporary one.
$rStart = $d.Range($r.Start, $r.Start)
$d.Bookmarks.Add("zzTempAnchor", $rStart)

$rEnd = $d.Range($r.End, $r.End)

# Adding a bookmark named "_GoBack" automatically replaces/removes any
# previously existing "_GoBack" bookmark elsewhere in the document,
# exactly mirroring Word's own behaviour of keeping only one instance of
# that bookmark, tracking the latest edit location.
$d.Bookmarks.Add("_GoBack", $rEnd)

$d.Bookmarks("zzTempAnchor").Delete()
